$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@('3535 Opal Meadow Heights Aged Care Community Meadow Heights', 25)
  ,@('7 Chefs Fawkner', 6)
  ,@('A1 Bakery Brunswick', 5)
  ,@('Acquire BPO Southbank', 9)
  ,@('Al Haj Halal Meats Glenroy', 50)
  ,@('Al-Taqwa College Truganina', 7)
  ,@('Amiga Montessori Craigieburn', 25)
  ,@('Best&Less Fountain Gate Narre Warren', 7)
  ,@('Budget Car and Truck Rentals Campbellfield', 7)
  ,@('CS Square Caroline Springs', 8)
  ,@('Campbellfield Heights Primary School Campbellfield', 6)
  ,@('Can Panel Cambellfield', 5)
  ,@('Cannie Road Construction Site Cannie', 7)
  ,@('Caroline Springs Police Station', 8)
  ,@('Cedars Medical Clinic Coburg', 44)
  ,@('Chemist Warehouse Fillo Drive Somerton', 5)
  ,@('City of Hobsons Bay Community', 9)
  ,@('City of Moreland Community', 7)
  ,@('City of Wyndham Community', 6)
  ,@('Classy Cabinets and Kitchens Craigieburn', 13)
  ,@('Coles Aurora Village Epping', 6)
  ,@('Coles Broadmeadows Central Shopping Centre', 10)
  ,@('Coles Campbellfield Plaza Campbellfield', 10)
  ,@('Coles Coburg North Village', 27)
  ,@('Coles Greenvale Shopping Centre', 7)
  ,@('Coles Pakenham Place Shopping Centre', 8)
  ,@('Coles Roxburgh Village Roxburgh Park', 11)
  ,@('Community Kids Meadow Heights', 19)
  ,@('Construction Site Olea Apartment Caulfield North', 8)
  ,@('Costco Wholesale Epping', 20)
  ,@('Croydon Orthodontics', 5)
  ,@('DRC Laverton Automotive Repairs Laverton North', 5)
  ,@('Direct Freight Express Cambellfield', 10)
  ,@('Elite Smart Community Care Campbellfield', 5)
  ,@('Fitzroy Community School Fitzroy North', 49)
  ,@('Glenroy West Primary School', 6)
  ,@('Hamilton Marino 236 Jasper Road McKinnon', 11)
  ,@('Health Care Providers Association South Melbourne', 13)
  ,@('IGA Meadow Heights Shopping Centre Meadow Heights', 7)
  ,@('Ilim College Glenroy Campus Hadfield', 22)
  ,@('Ilim Learning Sanctuary Glenroy', 15)
  ,@('Industrial Galvanizers Valmont Coatings Campbellfield', 17)
  ,@('Islamic College of Melbourne Tarneit', 8)
  ,@('KFC Fawkner', 9)
  ,@('Kasr Sweets Coolaroo', 9)
  ,@('Kids House Early Learning Cheltenham', 13)
  ,@('Learning Nest Early Learning Centre Meadow Heights', 9)
  ,@('Level Crossing Removal Project Lilydale Construction Site John Street', 8)
  ,@('Lineage Logistics Laverton North', 7)
  ,@('Linfox Somerton National Distribution Centre Somerton', 6)
  ,@('McDonalds Thomastown II', 7)
  ,@('Melbourne Metropolitan Remand Centre Ravenhall', 9)
  ,@('Melbourne Truck Repairs Campbellfield', 6)
  ,@('Melbourne West Police Station Docklands', 6)
  ,@('Melbourne Youth Justice Centre Parkville', 5)
  ,@('Melton Police Station Melton', 5)
  ,@('Mercy Hospital for Women Heidelberg', 6)
  ,@('Mernda YMCA Early Learning Centre Mernda', 5)
  ,@('Montessori Beginnings Greenvale', 6)
  ,@('MyCentre Childcare Broadmeadows', 15)
  ,@('National Gallery of Victoria Melbourne', 9)
  ,@('Newbury Child and Community Centre Craigieburn', 10)
  ,@('Newport Football Club Altona North', 5)
  ,@('Nido Early School Moonee Ponds', 18)
  ,@('Nino Early Learning Adventures Lalor', 5)
  ,@('North Geelong House Party', 7)
  ,@('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 52)
  ,@('OnQ Plumbing and Excavations Craigieburn', 15)
  ,@('Oporto Coolaroo', 16)
  ,@('Paisley Park Early Learning Centre Bundoora', 10)
  ,@('Panorama Construction Site Whitehorse Rd Box Hill', 27)
  ,@('People First Healthcare Home Residence Disability Support Taylors Lakes', 5)
  ,@('Ramsay Health Care Warringal Private Hospital Heidelberg', 9)
  ,@('Richmond Quarter 261-271 Bridge Road Construction Site Richmond', 9)
  ,@('Salta Drive Construction Site Rangedale Drainage Altona North', 6)
  ,@('St Vincents Hospital Emergency Department Melbourne', 6)
  ,@('Tek Foods Somerton', 11)
  ,@('The Homestead Child and Family Centre Roxburgh Park', 18)
  ,@('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 8)
  ,@('Tip Top Warehouse Dandenong', 8)
  ,@('Tunstall Fresh Tunstall Square Shopping Centre Doncaster East', 5)
  ,@('Unilodge College Square Student Accommodation 570 Lygon Street Carlton', 14)
  ,@('Werribee Mercy Hospital Emergency Department', 7)
  ,@('Western Health Footscray Hospital Emergency Department', 5)
  ,@('Western Health Sunshine Hospital Emergency Department', 7)
  ,@('Woodlands Long Day Care and Kindergarten Roxburgh Park', 5)
  ,@('Woolworths Greenvale Lakes Roxburgh Park', 6)
  ,@('Yara Childcare Centre Truganina', 13)
)

$n = $data.Count
$arr = New-Object 'object[,]' $n,2
for ($i = 0; $i -lt $n; $i++) {
  $arr[$i,0] = $data[$i][0]
  $arr[$i,1] = $data[$i][1]
}
$ws.Range("A2:B89").Value = $arr

$ws.Range("A90:B92").ClearContents() | Out-Null

Write-Host "Done. Dimensions:" $ws.Range("A1:B89").Address()